$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 11:22"

# Belgica (row 16): Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$ws.Range("B16").Value = 47859
$ws.Range("C16").Value = 525
$ws.Range("D16").Value = 11283
$ws.Range("E16").Value = 29075
$ws.Range("F16").Value = 797
$ws.Range("G16").Value = 170
$ws.Range("H16").Value = 7501

# Indonesia (row 40)
$ws.Range("B40").Value = 9771
$ws.Range("C40").Value = 260
$ws.Range("D40").Value = 1391
$ws.Range("E40").Value = 7596
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 11
$ws.Range("H40").Value = 784

# Estonia (row 76)
$ws.Range("B76").Value = 1666
$ws.Range("C76").Value = 6
$ws.Range("D76").Value = 236
$ws.Range("E76").Value = 1380
$ws.Range("F76").Value = 10
$ws.Range("G76").Value = 0
